# Update the '想去人数' (want-to-go count) figures (column F) across all four
# sheets to reflect refreshed data output, per commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 21285
$ws.Cells.Item(4, 6).Value = 105
$ws.Cells.Item(5, 6).Value = 3115
$ws.Cells.Item(6, 6).Value = 817
$ws.Cells.Item(7, 6).Value = 621
$ws.Cells.Item(8, 6).Value = 523
$ws.Cells.Item(9, 6).Value = 776
$ws.Cells.Item(10, 6).Value = 286
$ws.Cells.Item(13, 6).Value = 122
$ws.Cells.Item(14, 6).Value = 532
$ws.Cells.Item(16, 6).Value = 277
$ws.Cells.Item(17, 6).Value = 22
$ws.Cells.Item(18, 6).Value = 430
$ws.Cells.Item(19, 6).Value = 78
$ws.Cells.Item(20, 6).Value = 30
$ws.Cells.Item(22, 6).Value = 43

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 22
$ws.Cells.Item(5, 6).Value = 342
$ws.Cells.Item(14, 6).Value = 157

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 6140
$ws.Cells.Item(3, 6).Value = 706
$ws.Cells.Item(4, 6).Value = 706
$ws.Cells.Item(5, 6).Value = 1656
$ws.Cells.Item(6, 6).Value = 59

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 6140
$ws.Cells.Item(3, 6).Value = 706
$ws.Cells.Item(4, 6).Value = 706
$ws.Cells.Item(5, 6).Value = 1656
$ws.Cells.Item(7, 6).Value = 21285
$ws.Cells.Item(8, 6).Value = 22
$ws.Cells.Item(10, 6).Value = 105
$ws.Cells.Item(12, 6).Value = 342
$ws.Cells.Item(13, 6).Value = 3115
$ws.Cells.Item(14, 6).Value = 817
$ws.Cells.Item(16, 6).Value = 59
$ws.Cells.Item(17, 6).Value = 621
$ws.Cells.Item(18, 6).Value = 523
$ws.Cells.Item(19, 6).Value = 776
$ws.Cells.Item(20, 6).Value = 286
$ws.Cells.Item(26, 6).Value = 122
$ws.Cells.Item(29, 6).Value = 532
$ws.Cells.Item(33, 6).Value = 277
$ws.Cells.Item(34, 6).Value = 157
$ws.Cells.Item(35, 6).Value = 157
$ws.Cells.Item(36, 6).Value = 22
$ws.Cells.Item(37, 6).Value = 430
$ws.Cells.Item(39, 6).Value = 78
$ws.Cells.Item(40, 6).Value = 30
$ws.Cells.Item(44, 6).Value = 43
